# Pedidos.xlsx — remove the duplicate "quantidade 0" order line (row 105)
# from the filtered list on Planilha1. Deleting the row shifts every
# subsequent row up by one, shrinking the used range/autofilter/filter
# database from 384 to 383 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 105 (Remessa 80266246 / Material 10493-ARI-I /
# Quantidade 0) — everything below shifts up one row.
$ws.Rows.Item(105).Delete()

# Re-point the hidden _xlnm._FilterDatabase defined name at the new,
# one-row-shorter range.
$wb.Names.Item(1).RefersTo = "=Planilha1!`$A`$1:`$A`$383"

# The worksheet's AutoFilter header range also needs to shrink to match.
# Turn the filter off and reapply it on the resized range with the same
# "equals 80266246" criteria (as a discrete value list so it serializes
# back as <filters><filter val="80266246"/></filters>, matching the
# original filter's shape).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:A383").AutoFilter(1, @("80266246"), 7)

# Restore the cursor/selection position as left by the editor.
[void]$ws.Range("C393").Select()
